$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current id/name pairs from A2:B41 into memory
$startRow = 2
$endRow = 41
$n = $endRow - $startRow + 1

$ids = @()
$names = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $ids += $ws.Cells.Item($r, 1).Value2
    $names += $ws.Cells.Item($r, 2).Value2
}

# Sort ascending by id (column A) using a manual, type-safe bubble sort
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt $n - $i - 1; $j++) {
        if ($ids[$j] -gt $ids[$j + 1]) {
            $tmpId = $ids[$j]
            $ids[$j] = $ids[$j + 1]
            $ids[$j + 1] = $tmpId

            $tmpName = $names[$j]
            $names[$j] = $names[$j + 1]
            $names[$j + 1] = $tmpName
        }
    }
}

# Write sorted values back into the same range, leaving everything else untouched
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $ids[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
}
